$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 0.15

$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.37

$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0.03

$ws.Range("D6").Value = 98

$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0.14

$ws.Range("C11").Value = 49
$ws.Range("D11").Value = 25
$ws.Range("E11").Value = 0

$ws.Range("C12").Value = 0
$ws.Range("E12").Value = 1

$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 0.1

$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0.27

$ws.Range("C17").Value = 1
$ws.Range("E17").Value = 0.37

$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 1

$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 0

$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 0.05

$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 0.05

$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 0.18

$ws.Range("D25").Value = 3

$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = 0.02

$ws.Range("C28").Value = 0
$ws.Range("E28").Value = 0

$ws.Range("C32").Value = 0
$ws.Range("E32").Value = 0

$ws.Range("C33").Value = 7
$ws.Range("D33").Value = 0

$ws.Range("C34").Value = 6
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 0.01

$ws.Range("C35").Value = 9
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 0
